# =====================================================================
# Blog_stack_value_results.xlsx edit
# - Rename Sheet1 -> "Current", Sheet2 -> "RG table"
# - Add a new "Season Log" sheet (with a yellow tab) and make it active
# - Update the "Current" sheet with the latest stack data (new cutoff
#   of $3500 for cheap DK hitters drops the 3rd ("???") stack entirely)
# =====================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Rename the existing sheets
# ---------------------------------------------------------------
$current = $wb.Worksheets.Item("Sheet1")
$current.Name = "Current"

$rgTable = $wb.Worksheets.Item("Sheet2")
$rgTable.Name = "RG table"

# ---------------------------------------------------------------
# 2. Update "Current" sheet values
# ---------------------------------------------------------------

# -- Table headers (row 1) --
$current.Range("A1").Value = "Toronto Blue Jays hitters (FD, DK)"
$current.Range("F1").Value = "Houston Astros hitters (FD, DK)"
$current.Range("K1").Value = "???"

# -- Table 1 (A:D) - Toronto Blue Jays hitters --
$current.Range("A3").Value = "Drury"
$current.Range("B3").Value = 2200
$current.Range("C3").Value = 3

$current.Range("A4").Value = "McKinney"
$current.Range("B4").Value = 2500
$current.Range("C4").Value = 0

$current.Range("A5").Value = "Hernandez"
$current.Range("B5").Value = 2800
$current.Range("C5").Value = 6

$current.Range("A6").Value = "Smoak"
$current.Range("B6").Value = 3600
$current.Range("C6").Value = 6

# -- Table 2 (F:I) - Houston Astros hitters --
$current.Range("F3").Value = "Bregman"
$current.Range("G3").Value = 4200
$current.Range("H3").Value = 24.2

$current.Range("F4").Value = "Brantley"
$current.Range("G4").Value = 3700
$current.Range("H4").Value = 12.5

$current.Range("F5").Value = "Correa"
$current.Range("G5").Value = 3000
$current.Range("H5").Value = 12.5

$current.Range("F6").Value = "Gurriel"
$current.Range("G6").Value = 2700
$current.Range("H6").Value = 6.2

# -- Table 3 (K:N) - no longer qualifies under the new cutoff, so the
#    whole stack is cleared (keeps the L/N formatting, like the other
#    blank template rows further down the sheet) --
$current.Range("K3:K6").ClearContents()
$current.Range("L3:L6").ClearContents()
$current.Range("M3:M6").ClearContents()

# -- restore selection / active-cell state for the sheet --
$current.Range("F1:I1").Select()

# ---------------------------------------------------------------
# 3. Add the new "Season Log" sheet at the end
# ---------------------------------------------------------------
$log = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$log.Name = "Season Log"
$log.Tab.Color = 65535

# -- Headers --
$log.Range("A1").Value = "Date"
$log.Range("B1").Value = "Blog Site"
$log.Range("C1").Value = "Stack"
$log.Range("D1").Value = "Value"
$log.Range("E1").Value = "Result"
$log.Range("A1:E1").Font.Bold = $true

$log.Range("G1").Value = "Blog Site"
$log.Range("H1").Value = "Successes"
$log.Range("I1").Value = "Failures"
$log.Range("J1").Value = "Success Rate"
$log.Range("G1:J1").Font.Bold = $true
$log.Range("G1:J3").Borders.LineStyle = 1
$log.Range("G1:J3").Borders.Weight = 2

# -- Log rows --
$log.Range("A2").Value = 43552
$log.Range("B2").Value = "Draftshot"
$log.Range("C2").Value = "Chicago Cubs hitters (FD, DK)"
$log.Range("D2").Value = 7.36
$log.Range("E2").Value = "Success"

$log.Range("A3").Value = 43552
$log.Range("B3").Value = "Draftshot"
$log.Range("C3").Value = "New York Yankees hitters (FD)"
$log.Range("D3").Value = 5.18
$log.Range("E3").Value = "Success"

$log.Range("A4").Value = 43552
$log.Range("B4").Value = "Draftshot"
$log.Range("C4").Value = "Toronto Blue Jays hitters (FD)"
$log.Range("D4").Value = 0.28
$log.Range("E4").Value = "Failure"

$log.Range("A5").Value = 43553
$log.Range("B5").Value = "RG"
$log.Range("C5").Value = "Los Angeles Angels righties (FD, DK)"
$log.Range("D5").Value = 2.96
$log.Range("E5").Value = "Failure"

$log.Range("A6").Value = 43553
$log.Range("B6").Value = "RG"
$log.Range("C6").Value = "Boston Red Sox hitters (FD, DK)"
$log.Range("D6").Value = 1.99
$log.Range("E6").Value = "Failure"

$log.Range("A7").Value = 43553
$log.Range("B7").Value = "RG"
$log.Range("C7").Value = "Detroit Tigers hitters (FD, DK)"
$log.Range("D7").Value = 1.46
$log.Range("E7").Value = "Failure"

$log.Range("A8").Value = 43554
$log.Range("B8").Value = "Draftshot"
$log.Range("C8").Value = "Chicago Cubs hitters (FD, DK)"
$log.Range("D8").Value = 5.93
$log.Range("E8").Value = "Success"

$log.Range("A9").Value = 43554
$log.Range("B9").Value = "Draftshot"
$log.Range("C9").Value = "Boston Red Sox hitters (FD, DK)"
$log.Range("D9").Value = 2.07
$log.Range("E9").Value = "Failure"

$log.Range("A10").Value = 43556
$log.Range("B10").Value = "RG"
$log.Range("C10").Value = "Los Angeles Dodgers righties (FD, DK)"
$log.Range("D10").Value = 0.93
$log.Range("E10").Value = "Failure"

$log.Range("A11").Value = 43556
$log.Range("B11").Value = "RG"
$log.Range("C11").Value = "Houston Astros hitters (FD, DK)"
$log.Range("D11").Value = 3.22
$log.Range("E11").Value = "Success"

$log.Range("A12").Value = 43556
$log.Range("B12").Value = "RG"
$log.Range("C12").Value = "Toronto Blue Jays hitters (FD, DK)"
$log.Range("D12").Value = 3.82
$log.Range("E12").Value = "Success"

$log.Range("A13").Value = 43557
$log.Range("B13").Value = "RG"
$log.Range("C13").Value = "Toronto Blue Jays hitters (FD, DK)"
$log.Range("D13").Value = 1.35
$log.Range("E13").Value = "Failure"

$log.Range("A14").Value = 43557
$log.Range("B14").Value = "RG"
$log.Range("C14").Value = "Houston Astros hitters (FD, DK)"
$log.Range("D14").Value = 4.07
$log.Range("E14").Value = "Success"

# date-ish formatting for the log's date / blog-site columns
$log.Range("A2:A14").NumberFormat = "d-mmm"
$log.Range("B2:B14").NumberFormat = "m/d/yyyy"

# -- Summary table --
$log.Range("G2").Value = "Draftshot"
$log.Range("H2").Formula = '=COUNTIFS(B:B,"Draftshot",E:E,"Success")'
$log.Range("I2").Formula = '=COUNTIFS(B:B,"Draftshot",E:E,"Failure")'
$log.Range("J2").Formula = '=H2 / (H2+I2)'

$log.Range("G3").Value = "RG"
$log.Range("H3").Formula = '=COUNTIFS(B:B,"RG",E:E,"Success")'
$log.Range("I3").Formula = '=COUNTIFS(B:B,"RG",E:E,"Failure")'
$log.Range("J3").Formula = '=H3 / (H3+I3)'

$log.Range("J2:J3").NumberFormat = "0.0%"

# -- column widths (approximate the authored layout) --
$log.Columns.Item(1).ColumnWidth = 12.140625
$log.Columns.Item(2).ColumnWidth = 12.140625
$log.Columns.Item(3).ColumnWidth = 40.85546875
$log.Columns.Item(4).ColumnWidth = 6.85546875
$log.Columns.Item(5).ColumnWidth = 11.42578125
$log.Columns.Item(7).ColumnWidth = 9.28515625
$log.Columns.Item(8).ColumnWidth = 9.7109375
$log.Columns.Item(9).ColumnWidth = 8
$log.Columns.Item(10).ColumnWidth = 12.140625

# -- selection / active sheet --
$log.Range("H9").Select()
$log.Activate()
